$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The point_location_* columns (C,D) need to move in front of the
# point_number_* columns (A,B), i.e. final column order becomes:
# C, D, A, B (while keeping each pair's internal left/right order and
# all per-cell formatting / column widths intact). Do this with two
# real column cut+insert moves so widths/styles travel with the data
# instead of being re-derived (which would lose precision).

# Step 1: move column D in front of column A -> D, A, B, C
$ws.Columns.Item(4).Cut()
$ws.Columns.Item(1).Insert(-4161)

# Step 2: move (old) column C, now at position 4, in front -> C, D, A, B
$ws.Columns.Item(4).Cut()
$ws.Columns.Item(1).Insert(-4161)

# Update the active selection to match: whole column D selected, active cell D1.
$ws.Columns.Item(4).Select()
